$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览"  -- plain value updates in column F (no structural changes)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1028
$ws1.Range("F4").Value  = 248
$ws1.Range("F5").Value  = 30
$ws1.Range("F6").Value  = 463
$ws1.Range("F7").Value  = 751
$ws1.Range("F11").Value = 418
$ws1.Range("F13").Value = 87
$ws1.Range("F14").Value = 877
$ws1.Range("F16").Value = 2014
$ws1.Range("F17").Value = 504
$ws1.Range("F18").Value = 7880
$ws1.Range("F19").Value = 603
$ws1.Range("F21").Value = 66
$ws1.Range("F22").Value = 98
$ws1.Range("F23").Value = 21
$ws1.Range("F24").Value = 233

# ---------------------------------------------------------------------------
# Sheet 2: "演出" -- value updates + one new row inserted before row 14
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 524
$ws2.Range("F9").Value = 127

# Insert a new row before row 14, pushing the existing rows 14-20 down to 15-21
$ws2.Rows.Item(14).Insert()

# Carry column A's formatting (border / bold / alignment) down into the new row
$ws2.Range("A13").Copy()
$ws2.Range("A14").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "2024-10-25"
$ws2.Range("B14").Style = "Normal"
$ws2.Range("C14").Value = "广州·新生代流媒体小天后野田爱实 2024 巡演"
$ws2.Range("D14").Value = "南洲路158号2F SD Livehouse"
$ws2.Range("E14").Value = "2024.10.25 20:00-10.25 22:00"
$ws2.Range("F14").Value = 0
$ws2.Range("G14").Value = 280
$ws2.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=91823"
$ws2.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202409/oN7FyQ8v1725347758464.jpeg"

# Column A is a simple sequential row counter (row number - 1); restore it for
# every data row now that a row has been inserted.
for ($r = 2; $r -le 21; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet 3: "本地生活" -- plain value updates in column F (no structural changes)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5560
$ws3.Range("F3").Value = 406
$ws3.Range("F4").Value = 395

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型" -- value updates + the same new row inserted before row 39
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 5560
$ws4.Range("F4").Value  = 406
$ws4.Range("F5").Value  = 395
$ws4.Range("F6").Value  = 524
$ws4.Range("F7").Value  = 1028
$ws4.Range("F10").Value = 248
$ws4.Range("F11").Value = 30
$ws4.Range("F12").Value = 463
$ws4.Range("F13").Value = 751
$ws4.Range("F19").Value = 418
$ws4.Range("F22").Value = 87
$ws4.Range("F24").Value = 877
$ws4.Range("F26").Value = 127
$ws4.Range("F28").Value = 2014
$ws4.Range("F29").Value = 504
$ws4.Range("F30").Value = 7880
$ws4.Range("F33").Value = 603
$ws4.Range("F35").Value = 66
$ws4.Range("F36").Value = 98
$ws4.Range("F38").Value = 21
$ws4.Range("F39").Value = 233

# Insert a new row before row 39, pushing the existing rows 39-47 down to 40-48
$ws4.Rows.Item(39).Insert()

# Carry column A's formatting (border / bold / alignment) down into the new row
$ws4.Range("A38").Copy()
$ws4.Range("A39").PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false

$ws4.Range("B39").NumberFormat = "@"
$ws4.Range("B39").Value = "2024-10-25"
$ws4.Range("B39").Style = "Normal"
$ws4.Range("C39").Value = "广州·新生代流媒体小天后野田爱实 2024 巡演"
$ws4.Range("D39").Value = "南洲路158号2F SD Livehouse"
$ws4.Range("E39").Value = "2024.10.25 20:00-10.25 22:00"
$ws4.Range("F39").Value = 0
$ws4.Range("G39").Value = 280
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=91823"
$ws4.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202409/oN7FyQ8v1725347758464.jpeg"

# Column A is a simple sequential row counter (row number - 1); restore it for
# every data row now that a row has been inserted.
for ($r = 2; $r -le 48; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

Write-Output "done"
